# Weekly refresh of "Hortaliza, Vega Monumental Concepción - Alcachofa" price records.
# Updates the Fecha/Variedad/Volumen/Precio* figures on the existing weekly rows (2-21)
# and appends two new sampled rows (22-23) for the week, per the upstream source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise existing weekly observations (rows 2-21) ---
# Row 2
$ws.Range("D2").Value = 44350
$ws.Range("H2").Value = 'Argentina(o)'
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15600
$ws.Range("N2").Value = '$/caja 50 unidades'
$ws.Range("P2").Value = 312
$ws.Range("Q2").Value = 50

# Row 3
$ws.Range("D3").Value = 44350
$ws.Range("J3").Value = 40

# Row 4
$ws.Range("D4").Value = 44383
$ws.Range("H4").Value = 'Argentina(o)'
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 17400
$ws.Range("N4").Value = '$/caja 50 unidades'
$ws.Range("P4").Value = 348
$ws.Range("Q4").Value = 50

# Row 5
$ws.Range("D5").Value = 44364
$ws.Range("H5").Value = 'Argentina(o)'
$ws.Range("K5").Value = 19000
$ws.Range("L5").Value = 20000
$ws.Range("M5").Value = 19500
$ws.Range("N5").Value = '$/caja 50 unidades'
$ws.Range("P5").Value = 390
$ws.Range("Q5").Value = 50

# Row 6
$ws.Range("D6").Value = 44364
$ws.Range("H6").Value = 'Española'
$ws.Range("K6").Value = 19000
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 19500
$ws.Range("N6").Value = '$/caja 30 unidades'
$ws.Range("P6").Value = 650
$ws.Range("Q6").Value = 30

# Row 8
$ws.Range("D8").Value = 44358
$ws.Range("H8").Value = 'Argentina(o)'
$ws.Range("K8").Value = 18000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 19000
$ws.Range("N8").Value = '$/caja 50 unidades'
$ws.Range("P8").Value = 380
$ws.Range("Q8").Value = 50

# Row 9
$ws.Range("D9").Value = 44358
$ws.Range("H9").Value = 'Española'
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 19000
$ws.Range("N9").Value = '$/caja 30 unidades'
$ws.Range("P9").Value = 633
$ws.Range("Q9").Value = 30

# Row 10
$ws.Range("D10").Value = 44433
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 14000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 14500
$ws.Range("P10").Value = 290

# Row 11
$ws.Range("D11").Value = 44397
$ws.Range("H11").Value = 'Española'
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14500
$ws.Range("N11").Value = '$/caja 30 unidades'
$ws.Range("P11").Value = 483
$ws.Range("Q11").Value = 30

# Row 12
$ws.Range("D12").Value = 44335
$ws.Range("K12").Value = 17000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 17500
$ws.Range("P12").Value = 583

# Row 13
$ws.Range("D13").Value = 44426
$ws.Range("H13").Value = 'Madrigal'
$ws.Range("J13").Value = 50
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 12600
$ws.Range("N13").Value = '$/caja 40 unidades'
$ws.Range("P13").Value = 315
$ws.Range("Q13").Value = 40

# Row 14
$ws.Range("D14").Value = 44421
$ws.Range("H14").Value = 'Española'
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 14000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 14500
$ws.Range("N14").Value = '$/caja 30 unidades'
$ws.Range("P14").Value = 483
$ws.Range("Q14").Value = 30

# Row 15
$ws.Range("D15").Value = 44420
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14500
$ws.Range("P15").Value = 483

# Row 16
$ws.Range("D16").Value = 44441
$ws.Range("H16").Value = 'Española'
$ws.Range("K16").Value = 13000
$ws.Range("L16").Value = 14000
$ws.Range("M16").Value = 13500
$ws.Range("N16").Value = '$/caja 30 unidades'
$ws.Range("P16").Value = 450
$ws.Range("Q16").Value = 30

# Row 17
$ws.Range("D17").Value = 44342
$ws.Range("K17").Value = 17000
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = 17500
$ws.Range("P17").Value = 583

# Row 18
$ws.Range("D18").Value = 44342
$ws.Range("H18").Value = 'Madrigal'
$ws.Range("K18").Value = 15000
$ws.Range("L18").Value = 16000
$ws.Range("M18").Value = 15500
$ws.Range("N18").Value = '$/caja 40 unidades'
$ws.Range("P18").Value = 388
$ws.Range("Q18").Value = 40

# Row 19
$ws.Range("D19").Value = 44428

# Row 20
$ws.Range("D20").Value = 44442
$ws.Range("H20").Value = 'Española'
$ws.Range("K20").Value = 14500
$ws.Range("M20").Value = 14750
$ws.Range("N20").Value = '$/caja 30 unidades'
$ws.Range("P20").Value = 492
$ws.Range("Q20").Value = 30

# Row 21
$ws.Range("D21").Value = 44435
$ws.Range("H21").Value = 'Argentina(o)'
$ws.Range("K21").Value = 14000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 14500
$ws.Range("N21").Value = '$/caja 50 unidades'
$ws.Range("P21").Value = 290
$ws.Range("Q21").Value = 50

# --- Append new rows 22-23 ---
# Row 22
$ws.Range("A22").Value = 11
$ws.Range("B22").Value = 'Vega Monumental Concepción'
$ws.Range("C22").Value = 'Bíobío'
$ws.Range("D22").Value = 44376
$ws.Range("D22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = 100112013
$ws.Range("G22").Value = 'Alcachofa'
$ws.Range("H22").Value = 'Española'
$ws.Range("I22").Value = 'Primera'
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 19000
$ws.Range("L22").Value = 20000
$ws.Range("M22").Value = 19500
$ws.Range("N22").Value = '$/caja 30 unidades'
$ws.Range("O22").Value = 'Provincia de Limarí'
$ws.Range("P22").Value = 650
$ws.Range("Q22").Value = 30
$ws.Range("R22").Value = 'Hortaliza'

# Row 23
$ws.Range("A23").Value = 11
$ws.Range("B23").Value = 'Vega Monumental Concepción'
$ws.Range("C23").Value = 'Bíobío'
$ws.Range("D23").Value = 44399
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 100112013
$ws.Range("G23").Value = 'Alcachofa'
$ws.Range("H23").Value = 'Española'
$ws.Range("I23").Value = 'Primera'
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 14000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 14500
$ws.Range("N23").Value = '$/caja 30 unidades'
$ws.Range("O23").Value = 'Provincia de Limarí'
$ws.Range("P23").Value = 483
$ws.Range("Q23").Value = 30
$ws.Range("R23").Value = 'Hortaliza'

